$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-converted to a number by Excel;
# force the cell to Text format first so the value is stored as a string, matching
# the original inlineStr cell type.
$textCells = @(
    @{ Ref = "D4"; Value = '1.004' }
    @{ Ref = "D5"; Value = '1.003' }
    @{ Ref = "D6"; Value = '306.56' }
    @{ Ref = "D7"; Value = '0.3890' }
    @{ Ref = "D8"; Value = '0.3813' }
    @{ Ref = "D9"; Value = '1.004' }
    @{ Ref = "D10"; Value = '49.22' }
    @{ Ref = "D11"; Value = '1.330' }
    @{ Ref = "D12"; Value = '0.08358' }
    @{ Ref = "D13"; Value = '23.50' }
    @{ Ref = "D14"; Value = '7.015' }
    @{ Ref = "D15"; Value = '0.00001267' }
    @{ Ref = "D16"; Value = '7.394' }
    @{ Ref = "D18"; Value = '95.06' }
    @{ Ref = "D19"; Value = '0.06862' }
    @{ Ref = "D20"; Value = '20.57' }
    @{ Ref = "D21"; Value = '6.850' }
    @{ Ref = "D23"; Value = '13.44' }
    @{ Ref = "D25"; Value = '2.329' }
    @{ Ref = "D26"; Value = '2.653' }
    @{ Ref = "D27"; Value = '22.19' }
    @{ Ref = "D28"; Value = '157.19' }
    @{ Ref = "D29"; Value = '8.619' }
    @{ Ref = "D30"; Value = '139.55' }
    @{ Ref = "D31"; Value = '5.296' }
    @{ Ref = "D32"; Value = '2.413' }
    @{ Ref = "D34"; Value = '0.07930' }
    @{ Ref = "D35"; Value = '6.739' }
    @{ Ref = "D36"; Value = '0.02865' }
    @{ Ref = "D37"; Value = '0.2655' }
    @{ Ref = "D38"; Value = '0.9367' }
    @{ Ref = "D39"; Value = '0.09144' }
    @{ Ref = "D40"; Value = '1.434' }
    @{ Ref = "D41"; Value = '9.798' }
    @{ Ref = "D42"; Value = '0.7448' }
    @{ Ref = "D43"; Value = '12.87' }
    @{ Ref = "D44"; Value = '15.90' }
    @{ Ref = "D45"; Value = '0.6803' }
    @{ Ref = "D46"; Value = '2.438' }
    @{ Ref = "D47"; Value = '4.080' }
    @{ Ref = "D48"; Value = '1.002' }
    @{ Ref = "D49"; Value = '0.08318' }
    @{ Ref = "D50"; Value = '1.241' }
    @{ Ref = "D51"; Value = '131.16' }
)
foreach ($item in $textCells) {
    $rng = $ws.Range($item.Ref)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
}

# Cells whose new value is naturally text (not numeric-looking), assign directly.
$plainCells = @(
    @{ Ref = "D2"; Value = '24.019.70' }
    @{ Ref = "E2"; Value = '  -3.66%  ' }
    @{ Ref = "D3"; Value = '1.635.62' }
    @{ Ref = "E3"; Value = '  -3.31%  ' }
    @{ Ref = "E4"; Value = '  -0.49%  ' }
    @{ Ref = "E5"; Value = '  -0.60%  ' }
    @{ Ref = "E6"; Value = '  -2.91%  ' }
    @{ Ref = "E7"; Value = '  -1.97%  ' }
    @{ Ref = "E8"; Value = '  -4.40%  ' }
    @{ Ref = "E9"; Value = '  -0.34%  ' }
    @{ Ref = "E10"; Value = '  -5.92%  ' }
    @{ Ref = "E11"; Value = '  -8.18%  ' }
    @{ Ref = "E12"; Value = '  -4.19%  ' }
    @{ Ref = "E13"; Value = '  -7.86%  ' }
    @{ Ref = "E14"; Value = '  -4.82%  ' }
    @{ Ref = "E15"; Value = '  -5.23%  ' }
    @{ Ref = "E16"; Value = '  -5.86%  ' }
    @{ Ref = "D17"; Value = '1.648.03' }
    @{ Ref = "E17"; Value = '  -3.35%  ' }
    @{ Ref = "E18"; Value = '  +0.21%  ' }
    @{ Ref = "E19"; Value = '  -4.69%  ' }
    @{ Ref = "E20"; Value = '  +0.79%  ' }
    @{ Ref = "E21"; Value = '  -4.52%  ' }
    @{ Ref = "E22"; Value = '  -0.51%  ' }
    @{ Ref = "E23"; Value = '  -5.19%  ' }
    @{ Ref = "D24"; Value = '24.025.15' }
    @{ Ref = "E24"; Value = '  -3.60%  ' }
    @{ Ref = "E25"; Value = '  -2.23%  ' }
    @{ Ref = "E26"; Value = '  -6.98%  ' }
    @{ Ref = "E27"; Value = '  -4.64%  ' }
    @{ Ref = "E28"; Value = '  -2.96%  ' }
    @{ Ref = "E29"; Value = '  +7.03%  ' }
    @{ Ref = "E30"; Value = '  -6.69%  ' }
    @{ Ref = "E31"; Value = '  -13.02%  ' }
    @{ Ref = "E32"; Value = '  -9.10%  ' }
    @{ Ref = "D33"; Value = '1.817.12' }
    @{ Ref = "E33"; Value = '  -5.59%  ' }
    @{ Ref = "E34"; Value = '  -6.74%  ' }
    @{ Ref = "E35"; Value = '  -3.88%  ' }
    @{ Ref = "E36"; Value = '  -7.62%  ' }
    @{ Ref = "E37"; Value = '  -6.81%  ' }
    @{ Ref = "E38"; Value = '  -9.17%  ' }
    @{ Ref = "E39"; Value = '  -5.26%  ' }
    @{ Ref = "B40"; Value = 'TrustWalletToken' }
    @{ Ref = "C40"; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Ref = "E40"; Value = '  -2.71%  ' }
    @{ Ref = "B41"; Value = 'FraxShare' }
    @{ Ref = "C41"; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Ref = "E41"; Value = '  -9.37%  ' }
    @{ Ref = "E42"; Value = '  -7.55%  ' }
    @{ Ref = "E43"; Value = '  -7.43%  ' }
    @{ Ref = "E44"; Value = '  -6.02%  ' }
    @{ Ref = "E45"; Value = '  -6.22%  ' }
    @{ Ref = "E46"; Value = '  -7.01%  ' }
    @{ Ref = "E47"; Value = '  -3.42%  ' }
    @{ Ref = "E48"; Value = '  -0.69%  ' }
    @{ Ref = "E49"; Value = '  -6.68%  ' }
    @{ Ref = "E50"; Value = '  -10.24%  ' }
    @{ Ref = "E51"; Value = '  -5.40%  ' }
)
foreach ($item in $plainCells) {
    $ws.Range($item.Ref).Value = $item.Value
}
